# Regenerate save_data: update column G ("K" - strikeouts) to use the
# true K count instead of the old "Strike#" value, for the 2021
# quintana_jose save_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G ("K")
$newK = @{
    2  = 0
    3  = 2
    4  = 3
    5  = 1
    6  = 6
    7  = 2
    8  = 1
    10 = 2
    11 = 0
    12 = 0
    13 = 6
    14 = 0
    15 = 3
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 3
    21 = 3
    22 = 5
    23 = 4
    24 = 7
    25 = 6
    26 = 9
    27 = 2
    28 = 8
    29 = 4
    30 = 6
    31 = 7
    32 = 4
    33 = 6
    35 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
